# ProjectPlanner.xlsx update
# "Added documentation to every class. Program is considered done."
# Bumps ACTUAL DURATION (col F) and PERCENT COMPLETE (col G) for several
# tasks to reflect completed work, and nudges the final testing task's
# percent-complete upward as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- Documentation (row 5) ---------------------------------------------
$ws.Range("G5").Value = 0.6

# --- Image Processing (row 10, summary) + its sub-task (row 13) --------
# F10 previously carried a stray percentage-format override (cellXfs
# entry used only by this cell); clear that by pulling the plain
# formatting from a sibling cell (F13) before writing the new value so
# the now-unused style entry can drop out of cellXfs on save.
$ws.Range("F13").Copy() | Out-Null
$ws.Range("F10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F10").Value = 4

$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 1

# --- Comm PC/ESP32 (row 15, summary) + sub-task (row 17) ---------------
$ws.Range("F15").Value = 4
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 1

# --- Remote (row 23, summary) + sub-tasks (rows 24-27) ------------------
$ws.Range("F23").Value = 3

$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 1

$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 1

$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1

$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1

# --- Testing (row 28) ----------------------------------------------------
$ws.Range("G28").Value = 0.85

# --- Restore the on-screen selection to where the author left off -------
$ws.Range("F28").Select() | Out-Null
